# Fixed update to excel issue
# - Rename "Requested quantity" header on "Weekly Quantity" sheet to "Weekly_PO_Qty"
# - Rename "Requested quantity" header on "Monthly Trend" sheet to "Monthly_PO_Qty"
# - Add a new "PO Forecast" sheet (ds, PO_Forecast, yhat_lower, yhat_upper) with forecast data

$wb = $excel.ActiveWorkbook

$wsWeekly  = $wb.Worksheets.Item("Weekly Quantity")
$wsMonthly = $wb.Worksheets.Item("Monthly Trend")

# --- Header renames -------------------------------------------------------
$wsWeekly.Range("B1").Value  = "Weekly_PO_Qty"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- New "PO Forecast" sheet ----------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsForecast = $wb.Worksheets.Add($null, $lastSheet)
$wsForecast.Name = "PO Forecast"

# Copy the header style (bold / centered / thin border) from an existing
# header cell, and the date style from an existing date column, so the new
# sheet reuses the workbook's existing cell formats.
$wsWeekly.Range("B1").Copy()
$wsForecast.Range("A1:D1").PasteSpecial(-4122)  # xlPasteFormats

$wsWeekly.Range("A2:A11").Copy()
$wsForecast.Range("A2:A19").PasteSpecial(-4122) # xlPasteFormats

# Header row
$wsForecast.Range("A1").Value = "ds"
$wsForecast.Range("B1").Value = "PO_Forecast"
$wsForecast.Range("C1").Value = "yhat_lower"
$wsForecast.Range("D1").Value = "yhat_upper"

# Data rows
$rows = @(
    @(45137.99999999999, 13,  -77.84095816398604, 105.8628152758877),
    @(45144.99999999999, 15,  -74.83901691013355, 104.5728433802946),
    @(45158.99999999999, 19,  -74.89371793936591, 110.1458320137991),
    @(45165.99999999999, 21,  -72.96424173472849, 109.5789482670745),
    @(45172.99999999999, 23,  -66.15041180180701, 116.5777958317823),
    @(45466.99999999999, 108, 9.900706407383025,  201.7306004138291),
    @(45487.99999999999, 114, 15.65196986495716,  201.9597997029096),
    @(45501.99999999999, 118, 23.64146117427979,  209.8893595568858),
    @(45508.99999999999, 120, 28.02318653210956,  206.9247607552477),
    @(45529.99999999999, 126, 35.76787096848312,  215.7557259483727),
    @(45536.99999999999, 128, 30.01866262249658,  222.4486187696066),
    @(45543.99999999999, 130, 35.84410082222061,  222.8377147738892),
    @(45550.99999999999, 132, 43.34061650252514,  225.4078232542352),
    @(45557.99999999999, 135, 45.29370161967881,  223.3086020326581),
    @(45564.99999999999, 137, 47.9700434750296,   229.4354395778672),
    @(45571.99999999999, 139, 46.31519248058231,  230.0188412251592),
    @(45578.99999999999, 141, 45.23232604621796,  229.1500347561039),
    @(45585.99999999999, 143, 45.5037430853649,   235.113534759177)
)

$r = 2
foreach ($row in $rows) {
    $wsForecast.Cells.Item($r, 1).Value = $row[0]
    $wsForecast.Cells.Item($r, 2).Value = $row[1]
    $wsForecast.Cells.Item($r, 3).Value = $row[2]
    $wsForecast.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

Write-Output "PO Forecast sheet created with $($rows.Count) data rows."
